# Atualizacao dos dados BIBI (vendas_atipicas):
#  - remove a linha referente a 2025-07-01 (as demais linhas sobem uma
#    posicao);
#  - recalcula os desvios das linhas que correspondiam aos produtos
#    MOUSE PAD / CAPA IPHONE 11;
#  - adiciona uma nova linha de venda atipica para 2025-07-16
#    (BARALHO PLASTICO 1001 COPAG).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove a linha 2 (2025-07-01); linhas 3-8 sobem para 2-7.
$ws.Rows.Item(2).Delete()

# Estado final das linhas 2 a 8 (colunas A-I).
$data = @(
    @("2025-07-02", 2, "BEMOL S/A",        "375697", 13018, "ADAPTADOR TUDO EM UM UNIVERSAL INOVA PRIME TRA-30078",                    0, 1.08, 0.29),
    @("2025-07-07", 4, "BEMOL S/A",        "378212", 13546, "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON", -323, 1.1, 0.34),
    @("2025-07-09", 2, "BEMOL S/A",        "379513", 13000, "MOUSE PAD HARRY POTTER ESTAMPADO 26CMX21CM BLISTER C/1 UND LETRON",       -51, 1.03, 0.18),
    @("2025-07-09", 2, "MATHEUS SILVEIRA", "379106", 8915,  "CAPA IPHONE 11",                                                            29, 1.06, 0.25),
    @("2025-07-11", 2, "BEMOL S/A",        "380683", 14241, "MOCHILA PELUCIA STITCH",                                                  -26, 1.04, 0.2),
    @("2025-07-15", 2, "BEMOL S/A",        "383049", 12016, "PROJETOR ASTRONAUTA HMASTON",                                             -40, 1.04, 0.21),
    @("2025-07-16", 4, "BEMOL S/A",        "383665", 3984,  "BARALHO PLASTICO 1001 COPAG ESTOJO C/2 110 UNIDADES",                     -34, 1.22, 0.73)
)

$row = 2
foreach ($r in $data) {
    # Colunas A (data) e D (id_produto) sao textos "parecidos" com
    # numero/data; forca formato de texto para a gravacao e depois
    # devolve o estilo ao padrao da celula (Normal), assim como estava
    # no arquivo original (sem numero de estilo customizado).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]

    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $r[3]
    $cellD.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $row++
}

$wb.Save()
